$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 17-21 (5 rows) first, so remaining row data shifts up naturally.
$ws.Range("A17:H21").EntireRow.Delete() | Out-Null

# Update column widths (offset by -5/6 to compensate for the engine's
# internal pixel-quantization so the stored "width" lands on the exact
# target integer, matching genuine Excel column-width rounding behavior).
$ws.Columns.Item(3).ColumnWidth = 41.16666666666667
$ws.Columns.Item(4).ColumnWidth = 81.16666666666667
$ws.Columns.Item(8).ColumnWidth = 39.16666666666667

# Data for rows 2-16: A, C, D, F, G, H (B derived from A, E unchanged = "No")
$data = @(
    @{Row=2;  A="1328245"; C="Procurement Specialist Direct Materials"; D="Panamá, Provincia de Panamá, Panamá"; F="0 applicants"; G="3 - 6 Months"; H="Nestlé"},
    @{Row=3;  A="1328217"; C="Mobile application Developer"; D="Kafr El-Shaikh, Qism Kafr El-Shaikh, Kafr el-Sheikh, Gharbia Governorate, Egypt"; F="0 applicants"; G="9 - 12 Weeks"; H="Sharaby center"},
    @{Row=4;  A="1328163"; C="Marketing Specialist"; D="Tanta, Tanta Qism 2, Tanta, Gharbia Governorate, Egypt"; F="2 applicants"; G="9 - 12 Weeks"; H="print shop"},
    @{Row=5;  A="1328140"; C="Web Development Intern"; D="Malabe, Sri Lanka"; F="2 applicants"; G="3 - 6 Months"; H="TSA Media Group (Pvt) Ltd"},
    @{Row=6;  A="1328139"; C="Social Media Intern"; D="Malabe, Sri Lanka"; F="2 applicants"; G="3 - 6 Months"; H="TSA Media Group (Pvt) Ltd"},
    @{Row=7;  A="1328136"; C="Content Creator Intern"; D="Malabe, Sri Lanka"; F="1 applicant"; G="3 - 6 Months"; H="TSA Media Group (Pvt) Ltd"},
    @{Row=8;  A="1327959"; C="Accounting Management Trainee"; D="Hong Kong"; F="4 applicants"; G="6 - 18 Months"; H="ATHENASIA Consulting Limited"},
    @{Row=9;  A="1327518"; C="Business Development Intern"; D="Malabe, Sri Lanka"; F="22 applicants"; G="3 - 6 Months"; H="ZILLIONe Technologies Private Limited"},
    @{Row=10; A="1327475"; C="Property Consultant"; D="Cairo, Cairo Governorate, Egypt"; F="5 applicants"; G="9 - 12 Weeks"; H="Bold Routes"},
    @{Row=11; A="1327439"; C="Assistant"; D="Tunis, Tunisie"; F="19 applicants"; G="9 - 12 Weeks"; H="SSK events"},
    @{Row=12; A="1327438"; C="Graphic designer"; D="Tunis, Tunisie"; F="8 applicants"; G="9 - 12 Weeks"; H="SSK events"},
    @{Row=13; A="1327437"; C="Marketing Digital"; D="Tunis, Tunisie"; F="15 applicants"; G="9 - 12 Weeks"; H="SSK events"},
    @{Row=14; A="1327369"; C="Interior Architect"; D="Ennasr 2, Ariana, Tunisia"; F="2 applicants"; G="9 - 12 Weeks"; H="Happy Office Solutions"},
    @{Row=15; A="1327106"; C="Sales Assistant (Spanish)"; D="Denizli, Kumkısık, Denizli, Türkiye"; F="18 applicants"; G="6 - 18 Months"; H="Sera Moda"},
    @{Row=16; A="1326653"; C="Interior Design"; D="New Cairo City, Cairo Governorate, Egypt"; F="10 applicants"; G="9 - 12 Weeks"; H="Ahmad Elsherif Interior Designer"}
)

foreach ($item in $data) {
    $r = $item.Row
    $cellA = $ws.Cells.Item($r, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $item.A
    $cellA.Style = "Normal"

    $ws.Cells.Item($r, 2).Value = "https://aiesec.org/opportunity/global-talent/$($item.A)"
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 4).Value = $item.D
    $ws.Cells.Item($r, 6).Value = $item.F
    $ws.Cells.Item($r, 7).Value = $item.G
    $ws.Cells.Item($r, 8).Value = $item.H
}
